# "Generate Report for Handoff"
# A new source file (6ce311f6-437e-467c-a86a-aacdd1524fc0.md) was handed off
# for localization. Insert one new row for it, right above the trailing
# ".localization-config" row, on all three sheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$newMd       = "6ce311f6-437e-467c-a86a-aacdd1524fc0.md"
$newZhXlf    = "6ce311f6-437e-467c-a86a-aacdd1524fc0.978dfe207d0544e8408b73645b7a47a5bab8eeb6.zh-cn.xlf"
$newDeXlf    = "6ce311f6-437e-467c-a86a-aacdd1524fc0.978dfe207d0544e8408b73645b7a47a5bab8eeb6.de-de.xlf"
$zhHandoffDt = "2016-03-10 07:50:59"
$deHandoffDt = "2016-03-10 07:51:08"
$epoch       = "0001-01-01 00:00:00"

$mdUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/a9dc797029e8817a59d0a9e5ba4edc178fef1f35/e2e/" + $newMd
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d0f1dcd9f678be1b842eebd7a5d3cedfddd4baf7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/" + $newZhXlf
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f8dfddea56517d0a4a7932f3bd80368ac47fdfc7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/" + $newDeXlf

# ---------------------------------------------------------------------------
# Sheet "Overview": columns A=File Name, B=zh-cn, C=de-de
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Hyperlinks.Delete()
$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = $newMd
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"

# Only the freshly-inserted row needs its hyperlink style applied; A2 keeps
# the style it already had, and A4 inherited the old A3 style on the shift.
$ws.Range("A3").Style = "HyperLink"

$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl, $null, $null, "322627cf-f28b-4ad4-bcf4-d45a3baf76c2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl, $null, $null, $newMd) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f683652a6df31f9481d41ad09968c74ee440a636/.localization-config", $null, $null, ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn": columns A..I
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Hyperlinks.Delete()
$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = $newMd
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = $newZhXlf
$ws.Range("D3").Value = $zhHandoffDt
$ws.Range("G3").Value = $epoch
$ws.Range("H3").Value = "Include"

# Only the freshly-inserted row needs its hyperlink style applied; row 2
# keeps what it already had, and row 4 inherited the old row-3 style on shift.
$ws.Range("A3").Style = "HyperLink"
$ws.Range("C3").Style = "HyperLink"

$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl, $null, $null, "322627cf-f28b-4ad4-bcf4-d45a3baf76c2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d0f1dcd9f678be1b842eebd7a5d3cedfddd4baf7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/322627cf-f28b-4ad4-bcf4-d45a3baf76c2.2f1186afabb1ffe9f55f429b19fad9d79f07a323.zh-cn.xlf", $null, $null, "322627cf-f28b-4ad4-bcf4-d45a3baf76c2.2f1186afabb1ffe9f55f429b19fad9d79f07a323.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl, $null, $null, $newMd) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), $zhXlfUrl, $null, $null, $newZhXlf) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f683652a6df31f9481d41ad09968c74ee440a636/.localization-config", $null, $null, ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de": columns A..I
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Hyperlinks.Delete()
$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = $newMd
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = $newDeXlf
$ws.Range("D3").Value = $deHandoffDt
$ws.Range("G3").Value = $epoch
$ws.Range("H3").Value = "Include"

# Only the freshly-inserted row needs its hyperlink style applied; row 2
# keeps what it already had, and row 4 inherited the old row-3 style on shift.
$ws.Range("A3").Style = "HyperLink"
$ws.Range("C3").Style = "HyperLink"

$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl, $null, $null, "322627cf-f28b-4ad4-bcf4-d45a3baf76c2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f8dfddea56517d0a4a7932f3bd80368ac47fdfc7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/322627cf-f28b-4ad4-bcf4-d45a3baf76c2.2f1186afabb1ffe9f55f429b19fad9d79f07a323.de-de.xlf", $null, $null, "322627cf-f28b-4ad4-bcf4-d45a3baf76c2.2f1186afabb1ffe9f55f429b19fad9d79f07a323.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl, $null, $null, $newMd) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), $deXlfUrl, $null, $null, $newDeXlf) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f683652a6df31f9481d41ad09968c74ee440a636/.localization-config", $null, $null, ".localization-config") | Out-Null
